# feat: add 2022-Q3 data
#
# 1. "总计" (sheet1) gets a new top data row for 2022-Q3, pushing the
#    existing 2022-Q2 / 2022-Q1 rows down by one.
# 2. A brand-new "2022-Q3" worksheet is inserted right after "总计"
#    (so the tab order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1), holding
#    the per-fund breakdown table for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: update "总计" summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Make room: push old row 3 -> row 4, old row 2 -> row 3 (values + format)
$summary.Range("A3:D3").Copy($summary.Range("A4"))
$summary.Range("A2:D2").Copy($summary.Range("A3"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 11
$summary.Range("D2").Value = 3.61

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 4
$summary.Range("D3").Value = 3.75

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 3

# ---------------------------------------------------------------------
# Step 2: insert the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

# Re-resolve sheet references by name now that a sheet was inserted
# (stale handles captured before Worksheets.Add can point at the wrong
# sheet afterwards).
$template = $wb.Worksheets.Item("2022-Q2")

# Clone header row + first few data rows (values + styling) from the
# "2022-Q2" sheet, which uses the same table layout/formatting.
$template.Range("B1:H1").Copy($newSheet.Range("B1"))
$template.Range("A2:H5").Copy($newSheet.Range("A2"))
# Stretch the data-row styling down to cover all 11 data rows (2-12).
$newSheet.Range("A5:H5").Copy($newSheet.Range("A6:H12"))

$q3Data = @(
    @("0", "002251", "华夏军工安全灵活配置混合A", "42.95", "94.59", "4.64", "1.9929", "8"),
    @("1", "000001", "华夏成长混合", "32.40", "75.59", "2.56", "0.8294", "9"),
    @("2", "013566", "华夏军工安全灵活配置混合C", "11.35", "94.59", "4.64", "0.5266", "8"),
    @("3", "011160", "富国质量成长6个月持有期混合A", "3.70", "85.89", "2.83", "0.1047", "8"),
    @("4", "014663", "富国创新发展两年定期开放混合A", "2.24", "71.47", "4.14", "0.0927", "3"),
    @("5", "002210", "创金合信量化多因子股票A", "2.39", "91.71", "1.22", "0.0292", "9"),
    @("6", "014664", "富国创新发展两年定期开放混合C", "0.33", "71.47", "4.14", "0.0137", "3"),
    @("7", "003865", "创金合信量化多因子股票C", "0.75", "91.71", "1.22", "0.0092", "9"),
    @("8", "011161", "富国质量成长6个月持有期混合C", "0.14", "85.89", "2.83", "0.0040", "8"),
    @("9", "009514", "创金合信同顺创业板精选股票C", "0.16", "92.10", "2.05", "0.0033", "10"),
    @("10", "009513", "创金合信同顺创业板精选股票A", "0.09", "92.10", "2.05", "0.0018", "10")
)

$r = 2
foreach ($row in $q3Data) {
    $newSheet.Cells.Item($r, 1).Value = [double]$row[0]
    # Column B (fund code) has significant leading zeros ("002251") so it
    # must stay text too, same quote-prefix trick as columns D-G below.
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    # Columns D,E,F,G look numeric ("42.95", "94.59", ...) but the source
    # workbook stores them as plain text, so force text via a quote
    # prefix to stop Excel from re-interpreting them as numbers.
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 8).Value = [double]$row[7]
    $r = $r + 1
}
